$d = $word.ActiveDocument

# Insert an empty run (<w:r><w:t></w:t></w:r>) into each of the five
# paragraphs that gained one in the target revision. Using the Paragraph
# object's own Range (not a manually built Document.Range) is important:
# it correctly appends the new run as the paragraph's own content instead
# of bleeding into the following paragraph when the paragraph is empty.

$targets = @(5, 12, 35, 52, 53)

foreach ($idx in $targets) {
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    $r.InsertAfter("")
}
